# Update Sheets via scheduled runner
# Applies cached market-price / profit recompute values per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3746.6843
$ws.Range("J62").Value = 5256
$ws.Range("L62").Value = 5256
$ws.Range("N62").Value = -6504
$ws.Range("H65").Value = 3746.6843
$ws.Range("J65").Value = 5256
$ws.Range("L65").Value = 26280
$ws.Range("N65").Value = -32520
$ws.Range("H98").Value = 792.3684
$ws.Range("I98").Value = 867.61536
$ws.Range("J98").Value = 629.3333
$ws.Range("K98").Value = 867.61536
$ws.Range("L98").Value = 629.3333
$ws.Range("M98").Value = 630.38464
$ws.Range("N98").Value = -3625.3333
$ws.Range("H122").Value = 792.3684
$ws.Range("I122").Value = 867.61536
$ws.Range("J122").Value = 629.3333
$ws.Range("K122").Value = 2602.84608
$ws.Range("L122").Value = 1887.9999
$ws.Range("M122").Value = -152.8460800000003
$ws.Range("N122").Value = -6787.9999
$ws.Range("H129").Value = 176928.97
$ws.Range("J129").Value = 186730.58
$ws.Range("L129").Value = 560191.74
$ws.Range("N129").Value = -570191.74
$ws.Range("H137").Value = 32633.562
$ws.Range("I137").Value = 1371.591
$ws.Range("J137").Value = 101409.9
$ws.Range("K137").Value = 4114.772999999999
$ws.Range("L137").Value = 304229.7
$ws.Range("M137").Value = -1564.772999999999
$ws.Range("N137").Value = -309329.7
$ws.Range("H138").Value = 1628.1666
$ws.Range("I138").Value = 659.75757
$ws.Range("J138").Value = 2188.8245
$ws.Range("K138").Value = 1979.27271
$ws.Range("L138").Value = 6566.4735
$ws.Range("M138").Value = 3160.72729
$ws.Range("N138").Value = -16846.4735

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18325.574
$ws.Range("I32").Value = 20850.932
$ws.Range("J32").Value = 1770.4445
$ws.Range("K32").Value = 20850.932
$ws.Range("L32").Value = 1770.4445
$ws.Range("M32").Value = -20563.932
$ws.Range("N32").Value = -2344.4445
$ws.Range("H45").Value = 2873.9443
$ws.Range("I45").Value = 2679.7856
$ws.Range("K45").Value = 2679.7856
$ws.Range("M45").Value = -2302.7856
$ws.Range("H102").Value = 1444.4286
$ws.Range("I102").Value = 1105
$ws.Range("K102").Value = 1105
$ws.Range("M102").Value = 517
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045
$ws.Range("H119").Value = 18466.445
$ws.Range("J119").Value = 18466.445
$ws.Range("L119").Value = 18466.445
$ws.Range("N119").Value = -28142.445
$ws.Range("H122").Value = 2243.389
$ws.Range("I122").Value = 1577.6666
$ws.Range("J122").Value = 3574.8333
$ws.Range("K122").Value = 4732.9998
$ws.Range("L122").Value = 10724.4999
$ws.Range("M122").Value = -2282.9998
$ws.Range("N122").Value = -15624.4999
$ws.Range("H132").Value = 14863.128
$ws.Range("I132").Value = 2041.5714
$ws.Range("J132").Value = 47499.816
$ws.Range("K132").Value = 6124.7142
$ws.Range("L132").Value = 142499.448
$ws.Range("M132").Value = -3594.7142
$ws.Range("N132").Value = -147559.448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 430
$ws.Range("I11").Value = 287.5
$ws.Range("K11").Value = 287.5
$ws.Range("M11").Value = -147.5
$ws.Range("H134").Value = 44732.96
$ws.Range("I134").Value = 55506.25
$ws.Range("K134").Value = 166518.75
$ws.Range("M134").Value = -163983.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 878.25
$ws.Range("I134").Value = 810.5263
$ws.Range("K134").Value = 2431.5789
$ws.Range("M134").Value = 103.4211

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 35775716
$ws.Range("J37").Value = 35775716
$ws.Range("L37").Value = 107327148
$ws.Range("N37").Value = -107327372
$ws.Range("H76").Value = 4843.3335
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4843.3335
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 14530.0005
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -15296.0005
$ws.Range("H79").Value = 4843.3335
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4843.3335
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 14530.0005
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -17182.0005
$ws.Range("H103").Value = 2750.5715
$ws.Range("I103").Value = 752.6667
$ws.Range("J103").Value = 6346.8
$ws.Range("K103").Value = 2258.0001
$ws.Range("L103").Value = 19040.4
$ws.Range("M103").Value = -1379.0001
$ws.Range("N103").Value = -20798.4
$ws.Range("H129").Value = 238905.72
$ws.Range("I129").Value = 661.25
$ws.Range("J129").Value = 385517.7
$ws.Range("K129").Value = 1983.75
$ws.Range("L129").Value = 1156553.1
$ws.Range("M129").Value = 3016.25
$ws.Range("N129").Value = -1166553.1
$ws.Range("H131").Value = 775.85
$ws.Range("J131").Value = 799.07294
$ws.Range("L131").Value = 2397.21882
$ws.Range("N131").Value = -12477.21882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 69346.87
$ws.Range("I132").Value = 71698.8
$ws.Range("J132").Value = 64937
$ws.Range("K132").Value = 215096.4
$ws.Range("L132").Value = 194811
$ws.Range("M132").Value = -212566.4
$ws.Range("N132").Value = -199871

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5165.25
$ws.Range("I7").Value = 3188.889
$ws.Range("J7").Value = 6782.273
$ws.Range("K7").Value = 3188.889
$ws.Range("L7").Value = 6782.273
$ws.Range("M7").Value = -3076.889
$ws.Range("N7").Value = -7006.273
$ws.Range("H68").Value = 2516.5833
$ws.Range("J68").Value = 3033.1667
$ws.Range("L68").Value = 3033.1667
$ws.Range("N68").Value = -4531.1667
$ws.Range("H71").Value = 2516.5833
$ws.Range("J71").Value = 3033.1667
$ws.Range("L71").Value = 15165.8335
$ws.Range("N71").Value = -22653.8335
$ws.Range("H126").Value = 5165.25
$ws.Range("I126").Value = 3188.889
$ws.Range("J126").Value = 6782.273
$ws.Range("K126").Value = 9566.667000000001
$ws.Range("L126").Value = 20346.819
$ws.Range("M126").Value = -7096.667000000001
$ws.Range("N126").Value = -25286.819
$ws.Range("H136").Value = 28866.945
$ws.Range("I136").Value = 36707.145
$ws.Range("J136").Value = 1426.25
$ws.Range("K136").Value = 110121.435
$ws.Range("L136").Value = 4278.75
$ws.Range("M136").Value = -107571.435
$ws.Range("N136").Value = -9378.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1873.8422
$ws.Range("I122").Value = 1787.75
$ws.Range("J122").Value = 2333
$ws.Range("K122").Value = 5363.25
$ws.Range("L122").Value = 6999
$ws.Range("M122").Value = -2913.25
$ws.Range("N122").Value = -11899
$ws.Range("H126").Value = 1357.1428
$ws.Range("I126").Value = 1440
$ws.Range("K126").Value = 4320
$ws.Range("M126").Value = -1850
$ws.Range("H132").Value = 2041.8334
$ws.Range("I132").Value = 1676
$ws.Range("K132").Value = 5028
$ws.Range("M132").Value = -2498
$ws.Range("H136").Value = 62502600
$ws.Range("I136").Value = 76925470
$ws.Range("K136").Value = 230776410
$ws.Range("M136").Value = -230773860

